$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F header: United States
$ws.Range("F1").Value = "United States"

# Row 2 - "Payments from high-income countries to compensate low-income countries for climate damages"
$ws.Range("B2").Value = 0.508064516129032
$ws.Range("C2").Value = 0.513870541611625
$ws.Range("D2").Value = 0.642725598526704
$ws.Range("E2").Value = 0.498447204968944
$ws.Range("F2").Value = 0.410896708286039

# Row 3 - "High-income countries funding renewable energy in low-income countries"
$ws.Range("B3").Value = 0.6
$ws.Range("C3").Value = 0.652575957727873
$ws.Range("D3").Value = 0.692449355432781
$ws.Range("E3").Value = 0.647515527950311
$ws.Range("F3").Value = 0.52894438138479

# Row 4 - "High-income countries contributing $100 billion per year..."
$ws.Range("B4").Value = 0.541935483870968
$ws.Range("C4").Value = 0.597093791281374
$ws.Range("D4").Value = 0.622467771639042
$ws.Range("E4").Value = 0.545031055900621
$ws.Range("F4").Value = 0.445327279606508

# Row 5 - "Cancellation of low-income countries' public debt"
$ws.Range("B5").Value = 0.359677419354839
$ws.Range("C5").Value = 0.298546895640687
$ws.Range("D5").Value = 0.451197053406998
$ws.Range("E5").Value = 0.402173913043478
$ws.Range("F5").Value = 0.304578130911843

# Row 6 - "Democratise international institutions..."
$ws.Range("B6").Value = 0.435483870967742
$ws.Range("C6").Value = 0.416116248348745
$ws.Range("D6").Value = 0.532228360957643
$ws.Range("E6").Value = 0.414596273291925
$ws.Range("F6").Value = 0.335981838819523

# Row 7 - "Removing tariffs on imports from low-income countries"
$ws.Range("B7").Value = 0.401612903225806
$ws.Range("C7").Value = 0.498018494055482
$ws.Range("D7").Value = 0.530386740331492
$ws.Range("E7").Value = 0.571428571428571
$ws.Range("F7").Value = 0.386681800983731

# Row 8 - "A minimum wage in all countries at 50% of local median wage"
$ws.Range("B8").Value = 0.543548387096774
$ws.Range("C8").Value = 0.535006605019815
$ws.Range("D8").Value = 0.616942909760589
$ws.Range("E8").Value = 0.503105590062112
$ws.Range("F8").Value = 0.416193719258418

# Row 9 - "Fight tax evasion by creating a global financial register..."
$ws.Range("B9").Value = 0.72741935483871
$ws.Range("C9").Value = 0.684280052840158
$ws.Range("D9").Value = 0.74585635359116
$ws.Range("E9").Value = 0.698757763975155
$ws.Range("F9").Value = 0.437381763147938

# Row 10 - "A maximum wealth limit of $10 billion for each human"
$ws.Range("B10").Value = 0.401612903225806
$ws.Range("C10").Value = 0.467635402906209
$ws.Range("D10").Value = 0.441988950276243
$ws.Range("E10").Value = 0.498447204968944
$ws.Range("F10").Value = 0.342792281498297
